$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 previously held only the section header "grandes regiões e unidades da
# federação" (shared string) with no data. The correction removes that header
# row entirely; all subsequent rows (7-38, "norte" .. "distrito federal") shift
# up by one, and the former last row (38, "distrito federal") disappears from
# the sheet along with the now-unused shared string.
$ws.Rows.Item(6).Delete()
